# Bump the default font sizes defined on the Slide Master's text styles
# (Title style level 1, and Body style levels 1-5).
#
# Font.Size is expressed in points (matching PowerPoint's COM object model),
# e.g. 32 -> sz="3200" in the underlying OOXML.
#
# Note: the indexed "(...)" call syntax is used instead of ".Item(...)"
# because this host only resolves the TextStyles/Levels collections through
# that calling convention; also, the Body text style is addressed as
# TextStyles(3) here (Title = 1, Body = 3).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# Title style - level 1: 28pt -> 32pt
$master.TextStyles(1).Levels(1).Font.Size = 32

# Body style - levels 1-5
$master.TextStyles(3).Levels(1).Font.Size = 24
$master.TextStyles(3).Levels(2).Font.Size = 20
$master.TextStyles(3).Levels(3).Font.Size = 18
$master.TextStyles(3).Levels(4).Font.Size = 16
$master.TextStyles(3).Levels(5).Font.Size = 14
